# Fix the loop-variable label in the "Nested for Loop Example" slides:
# the printed/displayed variable name was "k" but the loop actually
# uses "i", so correct "k = " -> "i = " everywhere it appears.

$p = $ppt.ActivePresentation

# Slide 11: the R code snippet -> print(paste("k = " , i, "; j = ", j))
$s11 = $p.Slides.Item(11)
$sh11 = $s11.Shapes.Item(2)
$tr11 = $sh11.TextFrame.TextRange
$tr11.Paragraphs(3).Runs(1).Text = '    print(paste("i = " , i, "; j = ", j))'

# Slide 12: the 8 lines of console Output
$s12 = $p.Slides.Item(12)
$sh12 = $s12.Shapes.Item(2)
$tr12 = $sh12.TextFrame.TextRange
$tr12.Paragraphs(1).Runs(1).Text = '[1] "i =  1 ; j =  1"'
$tr12.Paragraphs(2).Runs(1).Text = '[1] "i =  1 ; j =  2"'
$tr12.Paragraphs(3).Runs(1).Text = '[1] "i =  2 ; j =  1"'
$tr12.Paragraphs(4).Runs(1).Text = '[1] "i =  2 ; j =  2"'
$tr12.Paragraphs(5).Runs(1).Text = '[1] "i =  3 ; j =  1"'
$tr12.Paragraphs(6).Runs(1).Text = '[1] "i =  3 ; j =  2"'
$tr12.Paragraphs(7).Runs(1).Text = '[1] "i =  4 ; j =  1"'
$tr12.Paragraphs(8).Runs(1).Text = '[1] "i =  4 ; j =  2"'
